# Insert a new daily price-report row above row 24 for
# "Vega Modelo de Temuco - Granada" (Fruta / hortaliza, semanal).
# Inserting the row shifts every existing row 24..123 down by one
# (to 25..124), which is exactly what the target diff shows: each
# row's D/K/L/M/N/O/P/Q/R/S/T values become the values previously
# held by the row above it, and a brand-new row 124 is created that
# carries what used to be row 123's data. Columns A,B,C,E,F,G,H,I,J
# are constant for every row in this sheet, so we just re-stamp them
# on the freshly inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(24).Insert()

$ws.Range("A24").Value = 10
$ws.Range("B24").Value = "Vega Modelo de Temuco"
$ws.Range("C24").Value = "La Araucanía"
$ws.Range("D24").Value = 44715
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100104
$ws.Range("H24").Value = "Frutos de pepita"
$ws.Range("I24").Value = 100104001
$ws.Range("J24").Value = "Granada"
$ws.Range("K24").Value = "Wonderfull"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 15000
$ws.Range("O24").Value = 16000
$ws.Range("P24").Value = 15400
$ws.Range("Q24").Value = "$/bandeja 15 kilos granel"
$ws.Range("R24").Value = "Provincia de Limarí"
$ws.Range("S24").Value = 1027
$ws.Range("T24").Value = 15
